$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold + bordered, same as H1) onto the
# two new header cells before writing their text, so I1/J1 match the
# look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# row, I-value, J-value for every data row (2..62)
$rows = @(
  @(2, 8, 8),
  @(3, 9, 9),
  @(4, 8, 9),
  @(5, 8, 8),
  @(6, 9, 9),
  @(7, 1, 1),
  @(8, 9, 9),
  @(9, 7, 7),
  @(10, 8, 8),
  @(11, 8, 9),
  @(12, 8, 8),
  @(13, 6, 6),
  @(14, 4, 4),
  @(15, 7, 7),
  @(16, 7, 8),
  @(17, 9, 9),
  @(18, 8, 8),
  @(19, 8, 8),
  @(20, 6, 7),
  @(21, 9, 9),
  @(22, 6, 8),
  @(23, 7, 7),
  @(24, 5, 5),
  @(25, 7, 8),
  @(26, 8, 8),
  @(27, 6, 6),
  @(28, 9, 9),
  @(29, 5, 5),
  @(30, 5, 5),
  @(31, 8, 8),
  @(32, 7, 7),
  @(33, 5, 6),
  @(34, 7, 7),
  @(35, 8, 8),
  @(36, 6, 6),
  @(37, 1, 2),
  @(38, 1, 2),
  @(39, 7, 7),
  @(40, 4, 6),
  @(41, 7, 7),
  @(42, 9, 9),
  @(43, 8, 9),
  @(44, 6, 7),
  @(45, 6, 7),
  @(46, 6, 6),
  @(47, 8, 8),
  @(48, 8, 8),
  @(49, 8, 9),
  @(50, 8, 8),
  @(51, 10, 10),
  @(52, 7, 7),
  @(53, 9, 9),
  @(54, 8, 8),
  @(55, 6, 7),
  @(56, 1, 1),
  @(57, 4, 4),
  @(58, 6, 6),
  @(59, 4, 4),
  @(60, 5, 6),
  @(61, 3, 3),
  @(62, 2, 2)
)

foreach ($r in $rows) {
  $ws.Cells.Item($r[0], 9).Value = $r[1]
  $ws.Cells.Item($r[0], 10).Value = $r[2]
}
